# Commit: Add PF/1.0.5 to meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 3) containing the new "PF/1.0.5" entry and the
# "X" markers for the sit2/uat2/prod columns.
$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
